$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 571031.7
$ws.Range("J17").Value = 571031.7
$ws.Range("L17").Value = 1713095.1
$ws.Range("N17").Value = -1713431.1
$ws.Range("H28").Value = 1135.8125
$ws.Range("I28").Value = 1258
$ws.Range("J28").Value = 769.25
$ws.Range("K28").Value = 1258
$ws.Range("L28").Value = 769.25
$ws.Range("M28").Value = -773
$ws.Range("N28").Value = -1739.25
$ws.Range("H62").Value = 96783.17999999999
$ws.Range("I62").Value = 114678.78
$ws.Range("J62").Value = 16253
$ws.Range("K62").Value = 114678.78
$ws.Range("L62").Value = 16253
$ws.Range("M62").Value = -114054.78
$ws.Range("N62").Value = -17501
$ws.Range("H65").Value = 96783.17999999999
$ws.Range("I65").Value = 114678.78
$ws.Range("J65").Value = 16253
$ws.Range("K65").Value = 573393.9
$ws.Range("L65").Value = 81265
$ws.Range("M65").Value = -570273.9
$ws.Range("N65").Value = -87505
$ws.Range("H107").Value = 691.94446
$ws.Range("I107").Value = 641
$ws.Range("J107").Value = 946.6667
$ws.Range("K107").Value = 641
$ws.Range("L107").Value = 946.6667
$ws.Range("M107").Value = 1279
$ws.Range("N107").Value = -4786.6667
$ws.Range("H113").Value = 4135.3
$ws.Range("I113").Value = 4458.8335
$ws.Range("J113").Value = 3650
$ws.Range("K113").Value = 4458.8335
$ws.Range("L113").Value = 3650
$ws.Range("M113").Value = -1204.8335
$ws.Range("N113").Value = -10158
$ws.Range("H121").Value = 882.3333
$ws.Range("J121").Value = 1199.6666
$ws.Range("L121").Value = 3598.9998
$ws.Range("N121").Value = -7092.9998
$ws.Range("H129").Value = 1278072.1
$ws.Range("I129").Value = 276.6
$ws.Range("J129").Value = 1950596.1
$ws.Range("K129").Value = 829.8000000000001
$ws.Range("L129").Value = 5851788.300000001
$ws.Range("M129").Value = 4170.2
$ws.Range("N129").Value = -5861788.300000001
$ws.Range("H132").Value = 2576.0908
$ws.Range("I132").Value = 3068.5908
$ws.Range("J132").Value = 1591.091
$ws.Range("K132").Value = 9205.7724
$ws.Range("L132").Value = 4773.272999999999
$ws.Range("M132").Value = -6675.7724
$ws.Range("N132").Value = -9833.272999999999
$ws.Range("H141").Value = 4073.2666
$ws.Range("I141").Value = 4166.5557
$ws.Range("J141").Value = 3933.3333
$ws.Range("K141").Value = 12499.6671
$ws.Range("L141").Value = 11799.9999
$ws.Range("M141").Value = -7319.667099999999
$ws.Range("N141").Value = -22159.9999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2405.1177
$ws.Range("I2").Value = 2505.375
$ws.Range("J2").Value = 2316
$ws.Range("K2").Value = 2505.375
$ws.Range("L2").Value = 2316
$ws.Range("M2").Value = -2392.375
$ws.Range("N2").Value = -2542
$ws.Range("H102").Value = 4503.75
$ws.Range("I102").Value = 4790
$ws.Range("J102").Value = 2500
$ws.Range("K102").Value = 4790
$ws.Range("L102").Value = 2500
$ws.Range("M102").Value = -3168
$ws.Range("N102").Value = -5744
$ws.Range("H110").Value = 790
$ws.Range("I110").Value = 790
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 790
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 1255
$ws.Range("N110").ClearContents()
$ws.Range("H116").Value = 2405.1177
$ws.Range("I116").Value = 2505.375
$ws.Range("J116").Value = 2316
$ws.Range("K116").Value = 2505.375
$ws.Range("L116").Value = 2316
$ws.Range("M116").Value = -211.375
$ws.Range("N116").Value = -6904
$ws.Range("H122").Value = 2172.25
$ws.Range("I122").Value = 883.375
$ws.Range("J122").Value = 4750
$ws.Range("K122").Value = 2650.125
$ws.Range("L122").Value = 14250
$ws.Range("M122").Value = -200.125
$ws.Range("N122").Value = -19150
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2405.1177
$ws.Range("I3").Value = 2505.375
$ws.Range("J3").Value = 2316
$ws.Range("K3").Value = 2505.375
$ws.Range("L3").Value = 2316
$ws.Range("M3").Value = -2391.375
$ws.Range("N3").Value = -2544
$ws.Range("H74").Value = 8500
$ws.Range("J74").Value = 8500
$ws.Range("L74").Value = 8500
$ws.Range("N74").Value = -10372
$ws.Range("H77").Value = 8500
$ws.Range("J77").Value = 8500
$ws.Range("L77").Value = 25500
$ws.Range("N77").Value = -34860
$ws.Range("H94").Value = 761.3182
$ws.Range("I94").Value = 842.1667
$ws.Range("K94").Value = 842.1667
$ws.Range("M94").Value = -391.1667
$ws.Range("H105").Value = 4838.5835
$ws.Range("I105").Value = 4006.45
$ws.Range("K105").Value = 4006.45
$ws.Range("M105").Value = -2259.45
$ws.Range("H134").Value = 21801.66
$ws.Range("I134").Value = 1684.8918
$ws.Range("J134").Value = 79057.08
$ws.Range("K134").Value = 5054.6754
$ws.Range("L134").Value = 237171.24
$ws.Range("M134").Value = -2519.6754
$ws.Range("N134").Value = -242241.24
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1642.8572
$ws.Range("I16").Value = 1120
$ws.Range("K16").Value = 1120
$ws.Range("M16").Value = -833
$ws.Range("H107").Value = 446.64706
$ws.Range("I107").Value = 553.2
$ws.Range("J107").Value = 402.25
$ws.Range("K107").Value = 553.2
$ws.Range("L107").Value = 402.25
$ws.Range("M107").Value = 1366.8
$ws.Range("N107").Value = -4242.25
$ws.Range("H113").Value = 1642.8572
$ws.Range("I113").Value = 1120
$ws.Range("K113").Value = 1120
$ws.Range("M113").Value = 1050
$ws.Range("H132").Value = 1637.0714
$ws.Range("I132").Value = 1171.0526
$ws.Range("J132").Value = 2620.889
$ws.Range("K132").Value = 3513.1578
$ws.Range("L132").Value = 7862.667
$ws.Range("M132").Value = -983.1578
$ws.Range("N132").Value = -12922.667
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 46.22222
$ws.Range("I12").Value = 57.5
$ws.Range("J12").Value = 37.2
$ws.Range("K12").Value = 172.5
$ws.Range("L12").Value = 111.6
$ws.Range("M12").Value = 0.5
$ws.Range("N12").Value = -457.6
$ws.Range("H122").Value = 633.087
$ws.Range("I122").Value = 520
$ws.Range("J122").Value = 646.87805
$ws.Range("K122").Value = 4680
$ws.Range("L122").Value = 5821.90245
$ws.Range("M122").Value = -2230
$ws.Range("N122").Value = -10721.90245
$ws.Range("H138").Value = 3061.4285
$ws.Range("I138").Value = 1612.5
$ws.Range("K138").Value = 4837.5
$ws.Range("M138").Value = 302.5
$ws.Range("H140").Value = 179704.88
$ws.Range("I140").Value = 217142.36
$ws.Range("K140").Value = 651427.08
$ws.Range("M140").Value = -646247.08
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 9990
$ws.Range("J33").Value = 9990
$ws.Range("L33").Value = 9990
$ws.Range("N33").Value = -10494
$ws.Range("H80").Value = 3735.182
$ws.Range("J80").Value = 3782.6667
$ws.Range("L80").Value = 3782.6667
$ws.Range("N80").Value = -5778.6667
$ws.Range("H83").Value = 3735.182
$ws.Range("J83").Value = 3782.6667
$ws.Range("L83").Value = 18913.3335
$ws.Range("N83").Value = -28897.3335
$ws.Range("H113").Value = 8377.143
$ws.Range("I113").Value = 12425
$ws.Range("J113").Value = 2980
$ws.Range("K113").Value = 12425
$ws.Range("L113").Value = 2980
$ws.Range("M113").Value = -10255
$ws.Range("N113").Value = -7320
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3368336.8
$ws.Range("I40").Value = 3368336.8
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 3368336.8
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -3368200.8
$ws.Range("N40").ClearContents()
$ws.Range("H61").Value = 3762.6858
$ws.Range("I61").Value = 4203.1333
$ws.Range("J61").Value = 1120
$ws.Range("K61").Value = 4203.1333
$ws.Range("L61").Value = 1120
$ws.Range("M61").Value = -4001.1333
$ws.Range("N61").Value = -1524
$ws.Range("H68").Value = 3350
$ws.Range("I68").Value = 3925
$ws.Range("J68").Value = 2200
$ws.Range("K68").Value = 3925
$ws.Range("L68").Value = 2200
$ws.Range("M68").Value = -3176
$ws.Range("N68").Value = -3698
$ws.Range("H71").Value = 3350
$ws.Range("I71").Value = 3925
$ws.Range("J71").Value = 2200
$ws.Range("K71").Value = 19625
$ws.Range("L71").Value = 11000
$ws.Range("M71").Value = -15881
$ws.Range("N71").Value = -18488
$ws.Range("H93").Value = 5100
$ws.Range("I93").Value = 5250
$ws.Range("J93").Value = 4800
$ws.Range("K93").Value = 5250
$ws.Range("L93").Value = 4800
$ws.Range("M93").Value = -4002
$ws.Range("N93").Value = -7296
$ws.Range("H113").Value = 3762.6858
$ws.Range("I113").Value = 4203.1333
$ws.Range("J113").Value = 1120
$ws.Range("K113").Value = 4203.1333
$ws.Range("L113").Value = 1120
$ws.Range("M113").Value = -2033.1333
$ws.Range("N113").Value = -5460
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 34994.43
$ws.Range("J69").Value = 34994.43
$ws.Range("L69").Value = 34994.43
$ws.Range("N69").Value = -36492.43
$ws.Range("H72").Value = 34994.43
$ws.Range("J72").Value = 34994.43
$ws.Range("L72").Value = 104983.29
$ws.Range("N72").Value = -112471.29
$ws.Range("H100").Value = 2483.3333
$ws.Range("I100").Value = 2483.3333
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 4966.6666
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -4425.6666
$ws.Range("N100").ClearContents()
$ws.Range("H122").Value = 1144.8182
$ws.Range("I122").Value = 1206.2858
$ws.Range("J122").Value = 1037.25
$ws.Range("K122").Value = 3618.8574
$ws.Range("L122").Value = 3111.75
$ws.Range("M122").Value = -1168.8574
$ws.Range("N122").Value = -8011.75
$ws.Range("H126").Value = 1160
$ws.Range("I126").Value = 1100
$ws.Range("J126").Value = 1175
$ws.Range("K126").Value = 3300
$ws.Range("L126").Value = 3525
$ws.Range("M126").Value = -830
$ws.Range("N126").Value = -8465
